# Report is now generated dynamically, so the previously hard-coded sample
# data (Forum Member / Number / Date / Amount rows) is cleared out, the
# rolled-over year label is bumped, and the now-unused detail rows in the
# second (prior-year) block are hidden so only the summary remains visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the hard-coded enrollment-fee detail rows for the first block
# (quarters laid out across columns A-C, E-G, I-K, M-O) ---
$ws.Range("A10:C19").ClearContents()
$ws.Range("E10:G19").ClearContents()
$ws.Range("I10:K19").ClearContents()
$ws.Range("M10:O19").ClearContents()

# --- Clear the hard-coded enrollment-fee detail rows further down the
# first block (only columns A-C and E-G still had sample data here) ---
$ws.Range("A31:C36").ClearContents()
$ws.Range("E31:G36").ClearContents()

# --- Bump the year label shown on the second (prior-year) block ---
$ws.Range("A72").Value = 2019

# --- The second block's detail rows are no longer needed now that the
# report builds its rows dynamically; hide them instead of deleting so the
# formulas/structure stay intact. Row 94 becomes a new, empty hidden row. ---
$ws.Rows("94:135").Hidden = $true
